$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates (rich-text shared strings) ---
# A8: "Volume 30   Number  7" -> "Volume 30   Number  9"
# (the trailing "7" is its own run; swap just that character)
$ws.Range("A8").Characters(21, 1).Text = "9"

# C9: "Report Covering the Week  2/13/2023  Through  2/19/2023"
#  -> "Report Covering the Week  2/27/2023  Through  3/5/2023"
$ws.Range("C9").Characters(27, 9).Text = "2/27/2023"
$ws.Range("C9").Characters(47, 9).Text = "3/5/2023"

# --- Weekly crime statistics table updates (rows 14-30) ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = -80
$ws.Range("F14").Value = 9
$ws.Range("G14").Value = 14
$ws.Range("H14").Value = -35.714285714285
$ws.Range("I14").Value = 18
$ws.Range("J14").Value = 22
$ws.Range("K14").Value = -18.181818181818
$ws.Range("L14").Value = -10
$ws.Range("M14").Value = 12.5
$ws.Range("N14").Value = -77.215189873417
$ws.Range("C15").Value = 4
$ws.Range("E15").Value = -42.857142857142
$ws.Range("F15").Value = 29
$ws.Range("G15").Value = 25
$ws.Range("H15").Value = 16
$ws.Range("I15").Value = 65
$ws.Range("J15").Value = 73
$ws.Range("K15").Value = -10.958904109589
$ws.Range("L15").Value = 12.068965517241
$ws.Range("M15").Value = 47.727272727272
$ws.Range("N15").Value = -33.673469387755
$ws.Range("C16").Value = 73
$ws.Range("D16").Value = 76
$ws.Range("E16").Value = -3.947368421052
$ws.Range("F16").Value = 283
$ws.Range("G16").Value = 347
$ws.Range("H16").Value = -18.443804034582
$ws.Range("I16").Value = 728
$ws.Range("J16").Value = 753
$ws.Range("K16").Value = -3.320053120849
$ws.Range("L16").Value = 28.395061728395
$ws.Range("M16").Value = 8.333333333333
$ws.Range("N16").Value = -74.678260869565
$ws.Range("C17").Value = 112
$ws.Range("D17").Value = 113
$ws.Range("E17").Value = -0.884955752212
$ws.Range("F17").Value = 495
$ws.Range("G17").Value = 477
$ws.Range("H17").Value = 3.773584905660
$ws.Range("I17").Value = 1150
$ws.Range("J17").Value = 1045
$ws.Range("K17").Value = 10.047846889952
$ws.Range("L17").Value = 26.651982378854
$ws.Range("M17").Value = 66.425470332850
$ws.Range("N17").Value = -10.921766072811
$ws.Range("C18").Value = 46
$ws.Range("D18").Value = 51
$ws.Range("E18").Value = -9.803921568627
$ws.Range("F18").Value = 221
$ws.Range("G18").Value = 240
$ws.Range("H18").Value = -7.916666666666
$ws.Range("I18").Value = 514
$ws.Range("J18").Value = 505
$ws.Range("K18").Value = 1.782178217821
$ws.Range("L18").Value = 46.857142857142
$ws.Range("M18").Value = -5.514705882352
$ws.Range("N18").Value = -84.022381100404
$ws.Range("C19").Value = 142
$ws.Range("D19").Value = 157
$ws.Range("E19").Value = -9.554140127388
$ws.Range("F19").Value = 532
$ws.Range("G19").Value = 613
$ws.Range("H19").Value = -13.213703099510
$ws.Range("I19").Value = 1216
$ws.Range("J19").Value = 1334
$ws.Range("K19").Value = -8.845577211394
$ws.Range("L19").Value = 29.087048832271
$ws.Range("M19").Value = 79.086892488954
$ws.Range("N19").Value = 4.647160068846
$ws.Range("C20").Value = 80
$ws.Range("D20").Value = 78
$ws.Range("E20").Value = 2.564102564102
$ws.Range("F20").Value = 393
$ws.Range("G20").Value = 349
$ws.Range("H20").Value = 12.607449856733
$ws.Range("I20").Value = 926
$ws.Range("J20").Value = 810
$ws.Range("K20").Value = 14.320987654321
$ws.Range("L20").Value = 167.630057803468
$ws.Range("M20").Value = 184.923076923077
$ws.Range("N20").Value = -65.703703703703
$ws.Range("C21").Value = 458
$ws.Range("D21").Value = 487
$ws.Range("E21").Value = -5.954825462012
$ws.Range("F21").Value = 1962
$ws.Range("G21").Value = 2065
$ws.Range("H21").Value = -4.987893462469
$ws.Range("I21").Value = 4617
$ws.Range("J21").Value = 4542
$ws.Range("K21").Value = 1.651254953764
$ws.Range("L21").Value = 44.68818552178
$ws.Range("M21").Value = 55.402221474251
$ws.Range("N21").Value = -59.578007354228
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 7
$ws.Range("E22").Value = -85.714285714285
$ws.Range("F22").Value = 25
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 41
$ws.Range("J22").Value = 65
$ws.Range("K22").Value = -36.923076923076
$ws.Range("L22").Value = 7.894736842105
$ws.Range("C23").Value = 33
$ws.Range("D23").Value = 28
$ws.Range("E23").Value = 17.857142857142
$ws.Range("F23").Value = 119
$ws.Range("G23").Value = 117
$ws.Range("H23").Value = 1.709401709401
$ws.Range("I23").Value = 289
$ws.Range("J23").Value = 247
$ws.Range("K23").Value = 17.004048582996
$ws.Range("L23").Value = 49.740932642487
$ws.Range("M23").Value = 79.503105590062
$ws.Range("C24").Value = 342
$ws.Range("D24").Value = 354
$ws.Range("E24").Value = -3.389830508474
$ws.Range("F24").Value = 1303
$ws.Range("G24").Value = 1396
$ws.Range("H24").Value = -6.661891117478
$ws.Range("I24").Value = 2876
$ws.Range("J24").Value = 2841
$ws.Range("K24").Value = 1.231960577261
$ws.Range("L24").Value = 32.290708371665
$ws.Range("M24").Value = 46.510443199184
$ws.Range("C25").Value = 201
$ws.Range("D25").Value = 207
$ws.Range("E25").Value = -2.898550724637
$ws.Range("F25").Value = 763
$ws.Range("G25").Value = 725
$ws.Range("H25").Value = 5.241379310344
$ws.Range("I25").Value = 1669
$ws.Range("J25").Value = 1568
$ws.Range("K25").Value = 6.441326530612
$ws.Range("L25").Value = 30.492572322126
$ws.Range("M25").Value = 4.377736085053
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -13.333333333333
$ws.Range("F26").Value = 49
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = -3.921568627450
$ws.Range("I26").Value = 112
$ws.Range("J26").Value = 126
$ws.Range("K26").Value = -11.111111111111
$ws.Range("L26").Value = 23.076923076923
$ws.Range("D27").Value = 19
$ws.Range("E27").Value = 21.052631578947
$ws.Range("F27").Value = 76
$ws.Range("G27").Value = 75
$ws.Range("H27").Value = 1.333333333333
$ws.Range("I27").Value = 185
$ws.Range("J27").Value = 142
$ws.Range("K27").Value = 30.281690140845
$ws.Range("L27").Value = 31.205673758865
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = -20
$ws.Range("F28").Value = 24
$ws.Range("G28").Value = 30
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 51
$ws.Range("J28").Value = 72
$ws.Range("K28").Value = -29.166666666666
$ws.Range("L28").Value = -12.068965517241
$ws.Range("M28").Value = -19.047619047619
$ws.Range("N28").Value = -72.131147540983
$ws.Range("C29").Value = 4
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -20
$ws.Range("F29").Value = 19
$ws.Range("G29").Value = 27
$ws.Range("H29").Value = -29.629629629629
$ws.Range("I29").Value = 41
$ws.Range("J29").Value = 67
$ws.Range("K29").Value = -38.805970149253
$ws.Range("L29").Value = -22.641509433962
$ws.Range("M29").Value = -24.074074074074
$ws.Range("N29").Value = -75
$ws.Range("J30").Value = 13
$ws.Range("K30").Value = -69.230769230769
$ws.Range("L30").Value = 0

# --- Row 30 (Hate Crimes): D/E change from numeric to text cells ---
# D30 becomes the text "0" (shared string already used by C30/F30)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("C30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null

# E30 becomes the text "***.*" (shared string already used by M30/N30)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("C30").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

Write-Host "Applied weekly crime data update."
